$d = $word.ActiveDocument

# Locate the paragraph holding the "!!!" marker (it is unique in the document).
$hit = $d.Content
$found = $hit.Find.Execute("!!!", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target marker text '!!!' in document"
}

$markerPara = $hit.Paragraphs(1)
$breakPara = $markerPara.Previous(1)

# Range spanning the stand-alone page-break paragraph and the "!!!" paragraph
# (start of the break paragraph through the end of the marker paragraph,
# including both paragraph marks).
$full = $d.Range($breakPara.Range.Start, $markerPara.Range.End)

# Replace both paragraphs with a single merged paragraph:
#  - takes the "Title" (a4) style + yellow-highlight rPr of the marker paragraph
#  - starts with the page-break run (no lastRenderedPageBreak on it anymore)
#  - then the highlighted "!!!" run (keeps its lastRenderedPageBreak)
#  - the _GoBack bookmark now sits at the end, right after the text
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a4"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:br w:type="page"/></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:lastRenderedPageBreak/><w:t>!!!</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$full.InsertXML($xml)
